# Product_backlog.xlsx update
# - Bumps the Sprint 1 "Day 2/3/4 Left" burndown numbers from 16 -> 17
# - Updates a status annotation in the task table (S -> COM)
# - Adds new "Building the environment" tasks (TUT2) for Diptin / Cristian / Haris / Bhuwan
#   into the previously-empty rows 39-44 of the Sprint 1 task list
# - Fills the helper numbering column (A) down through row 47

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# --- Burndown values: day 2/3/4 left (16 -> 17) ---
$ws.Range("B11").Value = 17
$ws.Range("B12").Value = 17
$ws.Range("B13").Value = 17

# --- Existing row 37 status annotation: S -> COM ---
$ws.Range("E37").Value = "COM"

# --- New task rows 39-42: Category=TUT2, Description="Building the environment" ---
$ws.Range("A39").Value = 22
$ws.Range("B39").Value = "TUT2"
$ws.Range("C39").Value = "Diptin"
$ws.Range("D39").Value = "Building the environment"
$ws.Range("E39").Value = "S"

$ws.Range("A40").Value = 23
$ws.Range("B40").Value = "TUT2"
$ws.Range("C40").Value = "Cristian"
$ws.Range("D40").Value = "Building the environment"

# Row 41 starts a fill-down sequence for the helper numbering column.
$ws.Range("A41:A47").Formula = "=A40+1"

$ws.Range("B41").Value = "TUT2"
$ws.Range("C41").Value = "Haris"
$ws.Range("D41").Value = "Building the environment"

$ws.Range("B42").Value = "TUT2"
$ws.Range("C42").Value = "Bhuwan"
$ws.Range("D42").Value = "Building the environment"

# Rows 45/46 got overwritten with plain (non-formula) numbers afterwards.
$ws.Range("A45").Value = 24
$ws.Range("A46").Value = 25

# Recalculate so dependent formulas (e.g. COUNT(A18:A45) and the chart caches) refresh.
$excel.Calculate()
